$d = $word.ActiveDocument

# --- Merge "Name: " + "Yuanji Jin " into a single run ---
$d.Content.Find.Execute("Name: ", $true, $false, $false, $false, $false, $true, 1, $false, "Name: ", 2) | Out-Null

# --- Merge "Student ID: " + "31942075" into a single run ---
$d.Content.Find.Execute("Student ID: ", $true, $false, $false, $false, $false, $true, 1, $false, "Student ID: ", 2) | Out-Null

# --- Merge "Studio: " + "20" into a single run ---
$d.Content.Find.Execute("Studio: ", $true, $false, $false, $false, $false, $true, 1, $false, "Studio: ", 2) | Out-Null

# --- Split "Vahid Pooryousef" into "Vahid " / proofErr-wrapped "Pooryousef",
#     and add a new paragraph with the course website link right after it ---
$tutorPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Tutor:*") {
        $tutorPara = $p
        break
    }
}

$tutorXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:r><w:t xml:space="preserve">Tutor: </w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">Vahid </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>Pooryousef</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '</w:p>' +
    '<w:p><w:r><w:t>https://jjdull.github.io/FIT3179/</w:t></w:r></w:p>'

$tutorPara.Range.InsertXML($tutorXml)

# --- Mark the run holding the picture as NoProof ---
$shape = $d.InlineShapes.Item(1)
$shape.Range.NoProofing = $true
